$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update required (minimum) hours per group/subject in column B
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 2
$ws.Range("B8").Value = 2
$ws.Range("B10").Value = 2
